$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 1000.0
$ws.Range("C4").Value = 103.0
$ws.Range("C6").Value = 28327.192654933424
$ws.Range("C7").Value = 27503.192654933424
$ws.Range("C8").Value = 27477.37687528542
$ws.Range("C9").Value = 9270.0
$ws.Range("C11").Value = 5928.407420745161
$ws.Range("C13").Value = 22398.785234188268
$ws.Range("C14").Value = 21574.785234188268
$ws.Range("C15").Value = 13128.785234188268
$ws.Range("C16").Value = 12822.727040188267
$ws.Range("C17").Value = 12047.197040188272
$ws.Range("C18").Value = 775.53
$ws.Range("C21").Value = 277794.86384950276
$ws.Range("C22").Value = 269714.1842495027
$ws.Range("C23").Value = 269461.0179340177
$ws.Range("C24").Value = 90907.64549999997
$ws.Range("C27").Value = 219657.0472168523
$ws.Range("C28").Value = 211576.36761685234
$ws.Range("C29").Value = 128749.40171685233
$ws.Range("C30").Value = 125747.99612866223
$ws.Range("C31").Value = 118142.6448541623
$ws.Range("C32").Value = 7605.351274499997

# --- FUSELAGE sheet ---
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C7").Value = 2882.0
$ws.Range("D7").Value = -11.662835249042121
$ws.Range("C8").Value = 3512.0
$ws.Range("D8").Value = 7.647509578544091
$ws.Range("C9").Value = 3744.0
$ws.Range("D9").Value = 14.758620689655203
$ws.Range("C12").Value = 3450.166666666666
$ws.Range("D12").Value = 5.75223499361432

# --- WING sheet ---
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value = 2695.0
$ws.Range("D7").Value = 23.908045977011543
$ws.Range("C13").Value = 2256.2857142857138
$ws.Range("D13").Value = 3.7372742200328477

# --- HORIZONTAL TAIL sheet ---
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C8").Value = 147.0
$ws.Range("D8").Value = -54.942528735632166
$ws.Range("C9").Value = 158.0
$ws.Range("D9").Value = -51.57088122605362
$ws.Range("C10").Value = 199.66666666666663
$ws.Range("D10").Value = -38.79948914431671

# --- LANDING GEARS sheet ---
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 909.0
$ws.Range("D5").Value = 4.482758620689682
$ws.Range("C6").Value = 1133.0
$ws.Range("D6").Value = 30.229885057471297
$ws.Range("C7").Value = 1279.0
$ws.Range("D7").Value = 47.011494252873604
$ws.Range("C8").Value = 1148.0
$ws.Range("D8").Value = 31.95402298850578
$ws.Range("C9").Value = 1117.25
$ws.Range("D9").Value = 28.419540229885065
